$d = $word.ActiveDocument

# Target fill color: srgbClr 56A12B (R=0x56=86, G=0xA1=161, B=0x2B=43)
# Word/VBA RGB() packs components as 0x00BBGGRR, i.e. R + (G*256) + (B*65536)
$newColor = 86 + (161 * 256) + (43 * 65536)

# Both "Rectangle 264" and "Rectangle 262" sticky-note header bars change
# their fill from #0dcc2b to #56a12b.
for ($i = 1; $i -le $d.Shapes.Count; $i++) {
    $shp = $d.Shapes.Item($i)
    if ($shp.Name -eq "Rectangle 264" -or $shp.Name -eq "Rectangle 262") {
        $shp.Fill.ForeColor.RGB = $newColor
    }
}
